# Complications documentation update
# Applies the set of paragraph insertions / lastRenderedPageBreak moves
# described by the commit "complications updated in documentation".

$d = $word.ActiveDocument
$wNS = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Find-ParagraphRange($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $text"
        return $null
    }
    $paraRng = $d.Range($rng.Start, $rng.End)
    $paraRng.Expand(4) | Out-Null   # wdParagraph
    return $paraRng
}

# 1. New paragraph "Airway injury" right before "Airway obstruction" (Chest section)
$target = Find-ParagraphRange("Airway obstruction")
$new = $target.InsertParagraphBefore()
$target2 = Find-ParagraphRange("Airway obstruction")
$prev = $target2.Previous(4, 1)
$prev.Text = "Airway injury"

# 2. New paragraph "Chyle leal" right after "Anastomotic leak- oesophagus" (Chest section)
$target = Find-ParagraphRange("Anastomotic leak- oesophagus")
$target.InsertParagraphAfter() | Out-Null
$target2 = Find-ParagraphRange("Anastomotic leak- oesophagus")
$nxt = $target2.Next(4, 1)
$nxt.Text = "Chyle leal"

# 3. Move lastRenderedPageBreak from "Deep venous thrombosis" to "Cellulitis" (Limbs section)
$target = Find-ParagraphRange("Deep venous thrombosis")
$xml = "<w:p xmlns:w='$wNS'><w:r><w:t>Deep venous thrombosis</w:t></w:r></w:p>"
$target.InsertXML($xml) | Out-Null

$target = Find-ParagraphRange("Cellulitis")
$xml = "<w:p xmlns:w='$wNS'><w:r><w:lastRenderedPageBreak/><w:t>Cellulitis</w:t></w:r></w:p>"
$target.InsertXML($xml) | Out-Null

# 4. New paragraph "Lymphoedema" right after "Joint infection" (Limbs section)
$target = Find-ParagraphRange("Joint infection")
$target.InsertParagraphAfter() | Out-Null
$target2 = Find-ParagraphRange("Joint infection")
$nxt = $target2.Next(4, 1)
$nxt.Text = "Lymphoedema"

# 5. New paragraph "Lymphocyst" (flagged by spellcheck) right before "Pancreatitis" (Abdomen and pelvis)
$target = Find-ParagraphRange("Pancreatitis")
$target.InsertParagraphBefore() | Out-Null
$target2 = Find-ParagraphRange("Pancreatitis")
$prev = $target2.Previous(4, 1)
$xml = "<w:p xmlns:w='$wNS'><w:proofErr w:type='spellStart'/><w:r><w:t>Lymphocyst</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$prev.InsertXML($xml) | Out-Null

# 6. Add lastRenderedPageBreak to "Peritonitis" (Abdomen and pelvis)
$target = Find-ParagraphRange("Peritonitis")
$xml = "<w:p xmlns:w='$wNS'><w:r><w:lastRenderedPageBreak/><w:t>Peritonitis</w:t></w:r></w:p>"
$target.InsertXML($xml) | Out-Null

# 7. Remove lastRenderedPageBreak from "Stoma ischemia" (Abdomen and pelvis)
$target = Find-ParagraphRange("Stoma ischemia")
$xml = "<w:p xmlns:w='$wNS'><w:r><w:t>Stoma ischemia</w:t></w:r></w:p>"
$target.InsertXML($xml) | Out-Null

# 8. New paragraphs "Nerve injury - facial" and "Nerve injury - recurrent laryngeal"
#    right after "Intracranial haematoma" (Head and Neck section)
$target = Find-ParagraphRange("Intracranial haematoma")
$target.InsertParagraphAfter() | Out-Null
$target2 = Find-ParagraphRange("Intracranial haematoma")
$nxt = $target2.Next(4, 1)
$nxt.Text = [char]0x2013
$nxt.Text = "Nerve injury " + [char]0x2013 + " facial"

$target3 = Find-ParagraphRange("Nerve injury " + [char]0x2013 + " facial")
$target3.InsertParagraphAfter() | Out-Null
$target4 = Find-ParagraphRange("Nerve injury " + [char]0x2013 + " facial")
$nxt2 = $target4.Next(4, 1)
$nxt2.Text = "Nerve injury " + [char]0x2013 + " recurrent laryngeal"

Write-Output "done"
